$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '44.432.97'
$ws.Range("E2").Value = '  +3.56%  '

# Row 3
$ws.Range("D3").Value = '2.282.13'
$ws.Range("E3").Value = '  +2.67%  '

# Row 4
$ws.Range("E4").Value = '  +0.17%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.25%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '106.28'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.47%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.592'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.10%  '

# Row 8
$ws.Range("E8").Value = '  +0.13%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.573'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.75%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.83'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.17%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0846'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.77%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.92'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.27%  '

# Row 13
$ws.Range("E13").Value = '  +1.34%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.888'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.53%  '

# Row 15
$ws.Range("D15").Value = '2.630.56'
$ws.Range("E15").Value = '  +2.94%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.67'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.56%  '

# Row 17
$ws.Range("D17").Value = '2.284.07'
$ws.Range("E17").Value = '  +2.45%  '

# Row 18
$ws.Range("D18").Value = '44.347.67'
$ws.Range("E18").Value = '  +3.72%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.05%  '

# Row 20
$ws.Range("E20").Value = '  +4.06%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.56'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.32%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.56'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.62%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.22'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.97%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '239.51'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.29%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.21'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.30%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.03%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.24'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.19%  '

# Row 28
$ws.Range("B28").Value = 'InjectiveProtocol'
$ws.Range("C28").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '38.80'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +13.13%  '

# Row 29
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.23'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.51%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.57'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.40%  '

# Row 31
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0892'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.25%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '162.87'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.95%  '

# Row 33
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.67'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.72%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.74'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.25%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.07'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.87%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.22'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.18%  '

# Row 37
$ws.Range("E37").Value = '  +12.24%  '

# Row 38
$ws.Range("E38").Value = '  -0.57%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.96'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.28%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.48'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.41%  '

# Row 41
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0330'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.15%  '

# Row 42
$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '15.57'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +24.86%  '

# Row 43
$ws.Range("E43").Value = '  +0.24%  '

# Row 44
$ws.Range("D44").Value = '1.777.57'
$ws.Range("E44").Value = '  -8.44%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.209'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.47%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '87.06'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.84%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.52'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.40%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '60.61'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.11%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.72'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.35%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '74.69'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.98%  '

# Row 51
$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.72'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.08%  '
